# Disponibilidad.xlsx - "Actualizar 02-05-2021 20-28-43" automated update
#
# 1) The previous check run (rows 772-785, timestamp 44232.83168540164) gets its
#    timestamp nudged very slightly to 44232.83168540509.
# 2) A brand-new check run is appended as rows 786-799 (same 14-service pattern)
#    stamped 44232.85306001319, each with its own hyperlink in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Tiny timestamp correction on the last existing block (rows 772-785)
# ---------------------------------------------------------------------------
$ws.Range("D772:D785").Value = 44232.83168540509

# ---------------------------------------------------------------------------
# 2) Append the new block: rows 786-799
# ---------------------------------------------------------------------------
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$addresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$subAddresses = @("","","","","","","","","/","","","","","")

$startRow = 786
$newTimestamp = 44232.85306001319

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i

    $ws.Range("A$row").Value = $names[$i]
    $ws.Range("B$row").Value = $urls[$i]
    $ws.Range("C$row").Value = "Disponible"
    $ws.Range("D$row").Value = $newTimestamp

    if ($subAddresses[$i] -ne "") {
        $ws.Hyperlinks.Add($ws.Range("B$row"), $addresses[$i], $subAddresses[$i])
    } else {
        $ws.Hyperlinks.Add($ws.Range("B$row"), $addresses[$i])
    }
}

# Re-apply the workbook's built-in Hyperlink cell style and the shared date
# number format so the new cells match the style of every other row instead
# of whatever ad-hoc style Hyperlinks.Add auto-generated.
$ws.Range("B786:B799").Style = "Hyperlink"
$ws.Range("D786:D799").NumberFormat = "YYYY-MM-DD HH:MM:SS"
